$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3640.1428
$ws.Range("I125").Value = 727
$ws.Range("J125").Value = 5825
$ws.Range("K125").Value = 6543
$ws.Range("L125").Value = 52425
$ws.Range("M125").Value = -4083
$ws.Range("N125").Value = -57345
$ws.Range("H132").Value = 2481.5
$ws.Range("I132").Value = 2505.138
$ws.Range("K132").Value = 7515.414
$ws.Range("M132").Value = -4985.414
$ws.Range("H140").Value = 416666.34
$ws.Range("J140").Value = 416666.34
$ws.Range("L140").Value = 416666.34
$ws.Range("N140").Value = -427026.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2970
$ws.Range("I45").Value = 1327.75
$ws.Range("K45").Value = 1327.75
$ws.Range("M45").Value = -950.75
$ws.Range("H61").Value = 2612.4736
$ws.Range("I61").Value = 2305.8
$ws.Range("K61").Value = 2305.8
$ws.Range("M61").Value = -2093.8
$ws.Range("H97").Value = 581.5333000000001
$ws.Range("I97").Value = 586.4286
$ws.Range("K97").Value = 586.4286
$ws.Range("M97").Value = -90.42859999999996
$ws.Range("H109").Value = 67845
$ws.Range("J109").Value = 67845
$ws.Range("L109").Value = 67845
$ws.Range("N109").Value = -70619
$ws.Range("H136").Value = 2612.4736
$ws.Range("I136").Value = 2305.8
$ws.Range("K136").Value = 6917.400000000001
$ws.Range("M136").Value = -4367.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2705.2415
$ws.Range("J105").Value = 3326.4546
$ws.Range("L105").Value = 3326.4546
$ws.Range("N105").Value = -6820.4546
$ws.Range("H134").Value = 2980070.8
$ws.Range("I134").Value = 3761936.8
$ws.Range("J134").Value = 8980
$ws.Range("K134").Value = 11285810.4
$ws.Range("L134").Value = 26940
$ws.Range("M134").Value = -11283275.4
$ws.Range("N134").Value = -32010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1315.1666
$ws.Range("I16").Value = 1378.8
$ws.Range("J16").Value = 997
$ws.Range("K16").Value = 1378.8
$ws.Range("L16").Value = 997
$ws.Range("M16").Value = -1091.8
$ws.Range("N16").Value = -1571
$ws.Range("H31").Value = 5578.4873
$ws.Range("I31").Value = 2471.2307
$ws.Range("J31").Value = 7132.115
$ws.Range("K31").Value = 2471.2307
$ws.Range("L31").Value = 7132.115
$ws.Range("M31").Value = -2176.2307
$ws.Range("N31").Value = -7722.115
$ws.Range("H34").Value = 5578.4873
$ws.Range("I34").Value = 2471.2307
$ws.Range("J34").Value = 7132.115
$ws.Range("K34").Value = 2471.2307
$ws.Range("L34").Value = 7132.115
$ws.Range("M34").Value = -2269.2307
$ws.Range("N34").Value = -7536.115
$ws.Range("H58").Value = 2852.2683
$ws.Range("I58").Value = 2671.9143
$ws.Range("K58").Value = 2671.9143
$ws.Range("M58").Value = -2468.9143
$ws.Range("H99").Value = 1884.2
$ws.Range("I99").Value = 1884.2
$ws.Range("K99").Value = 1884.2
$ws.Range("M99").Value = -386.2
$ws.Range("H113").Value = 1315.1666
$ws.Range("I113").Value = 1378.8
$ws.Range("J113").Value = 997
$ws.Range("K113").Value = 1378.8
$ws.Range("L113").Value = 997
$ws.Range("M113").Value = 791.2
$ws.Range("N113").Value = -5337
$ws.Range("H126").Value = 1884.2
$ws.Range("I126").Value = 1884.2
$ws.Range("K126").Value = 5652.6
$ws.Range("M126").Value = -3182.6
$ws.Range("H132").Value = 3734.52
$ws.Range("I132").Value = 3678.45
$ws.Range("J132").Value = 3958.8
$ws.Range("K132").Value = 11035.35
$ws.Range("L132").Value = 11876.4
$ws.Range("M132").Value = -8505.349999999999
$ws.Range("N132").Value = -16936.4
$ws.Range("H134").Value = 2811.6365
$ws.Range("I134").Value = 2592.8
$ws.Range("K134").Value = 7778.400000000001
$ws.Range("M134").Value = -5243.400000000001
$ws.Range("H136").Value = 2852.2683
$ws.Range("I136").Value = 2671.9143
$ws.Range("K136").Value = 8015.742899999999
$ws.Range("M136").Value = -5465.742899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 862.6
$ws.Range("H39").Value = 5003
$ws.Range("J39").Value = 5003
$ws.Range("L39").Value = 15009
$ws.Range("N39").Value = -15597
$ws.Range("H132").Value = 2002509.6
$ws.Range("J132").Value = 5004999
$ws.Range("L132").Value = 45044991
$ws.Range("N132").Value = -45050051
$ws.Range("H136").Value = 2142.5
$ws.Range("I136").Value = 2142.5
$ws.Range("K136").Value = 6427.5
$ws.Range("M136").Value = -1327.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2854.7
$ws.Range("I102").Value = 2854.7
$ws.Range("K102").Value = 2854.7
$ws.Range("M102").Value = -1232.7
$ws.Range("H113").Value = 12962.223
$ws.Range("I113").Value = 2035.8572
$ws.Range("J113").Value = 51204.5
$ws.Range("K113").Value = 2035.8572
$ws.Range("L113").Value = 51204.5
$ws.Range("M113").Value = 134.1428000000001
$ws.Range("N113").Value = -55544.5
$ws.Range("H122").Value = 1608.75
$ws.Range("I122").Value = 1421.25
$ws.Range("K122").Value = 4263.75
$ws.Range("M122").Value = -1813.75
$ws.Range("H126").Value = 2466.111
$ws.Range("I126").Value = 2099.2856
$ws.Range("K126").Value = 6297.8568
$ws.Range("M126").Value = -3827.8568
$ws.Range("H132").Value = 2472.5757
$ws.Range("I132").Value = 2136.5
$ws.Range("J132").Value = 5833.3335
$ws.Range("K132").Value = 6409.5
$ws.Range("L132").Value = 17500.0005
$ws.Range("M132").Value = -3879.5
$ws.Range("N132").Value = -22560.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1643.625
$ws.Range("I22").Value = 637.5
$ws.Range("J22").Value = 2649.75
$ws.Range("K22").Value = 637.5
$ws.Range("L22").Value = 2649.75
$ws.Range("M22").Value = -342.5
$ws.Range("N22").Value = -3239.75
$ws.Range("H27").Value = 1643.625
$ws.Range("I27").Value = 637.5
$ws.Range("J27").Value = 2649.75
$ws.Range("K27").Value = 637.5
$ws.Range("L27").Value = 2649.75
$ws.Range("M27").Value = -530.5
$ws.Range("N27").Value = -2863.75
$ws.Range("H40").Value = 33338034
$ws.Range("I40").Value = 55558640
$ws.Range("J40").Value = 7125
$ws.Range("K40").Value = 55558640
$ws.Range("L40").Value = 7125
$ws.Range("M40").Value = -55558504
$ws.Range("N40").Value = -7397
$ws.Range("H46").Value = 3235.9167
$ws.Range("I46").Value = 749.3333
$ws.Range("J46").Value = 3591.1428
$ws.Range("K46").Value = 749.3333
$ws.Range("L46").Value = 3591.1428
$ws.Range("M46").Value = -561.3333
$ws.Range("N46").Value = -3967.1428
$ws.Range("H61").Value = 1839.7
$ws.Range("I61").Value = 1921.7778
$ws.Range("J61").Value = 1772.5454
$ws.Range("K61").Value = 1921.7778
$ws.Range("L61").Value = 1772.5454
$ws.Range("M61").Value = -1719.7778
$ws.Range("N61").Value = -2176.5454
$ws.Range("H95").Value = 106000
$ws.Range("J95").Value = 106000
$ws.Range("L95").Value = 106000
$ws.Range("N95").Value = -111492
$ws.Range("H113").Value = 1839.7
$ws.Range("I113").Value = 1921.7778
$ws.Range("J113").Value = 1772.5454
$ws.Range("K113").Value = 1921.7778
$ws.Range("L113").Value = 1772.5454
$ws.Range("M113").Value = 248.2221999999999
$ws.Range("N113").Value = -6112.5454
$ws.Range("H121").Value = 82973.5
$ws.Range("J121").Value = 82973.5
$ws.Range("L121").Value = 82973.5
$ws.Range("N121").Value = -86467.5
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920
$ws.Range("H132").Value = 3508
$ws.Range("I132").Value = 3011
$ws.Range("K132").Value = 9033
$ws.Range("M132").Value = -6503

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 55000
$ws.Range("J97").Value = 55000
$ws.Range("L97").Value = 55000
$ws.Range("N97").Value = -56982
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2838.5715
$ws.Range("I132").Value = 2794.0908
$ws.Range("K132").Value = 8382.2724
$ws.Range("M132").Value = -5852.2724
